$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected; insert it as a new row at
# position 81 (pushing the existing rows 81-122 down to 82-123).
$ws.Rows.Item(81).Insert()

$ws.Range("A81").Value = 11
$ws.Range("B81").Value = "Vega Monumental Concepción"
$ws.Range("C81").Value = "Bíobío"
$ws.Range("D81").Value = Get-Date -Year 2022 -Month 8 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = 100112021
$ws.Range("G81").Value = "Ají"
$ws.Range("H81").Value = "Inferno"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 22
$ws.Range("K81").Value = 12000
$ws.Range("L81").Value = 14000
$ws.Range("M81").Value = 12909
$ws.Range("N81").Value = "$/caja 12 kilos"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 1076
$ws.Range("Q81").Value = 12
$ws.Range("R81").Value = "Hortaliza"
